$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  "SpiceJet",  "21:40", "02 h 45 m", "₹ 8,982"),
    @(3,  "SpiceJet",  "08:30", "05 h 30 m", "₹ 8,982"),
    @(4,  "Air India", "06:00", "02 h 55 m", "₹ 8,983"),
    @(5,  "IndiGo",    "06:20", "02 h 50 m", "₹ 8,983"),
    @(6,  "Vistara",   "07:05", "02 h 50 m", "₹ 8,983"),
    @(7,  "IndiGo",    "08:45", "02 h 45 m", "₹ 8,983"),
    @(8,  "Air India", "09:55", "03 h 05 m", "₹ 8,983"),
    @(9,  "Vistara",   "10:35", "02 h 50 m", "₹ 8,983"),
    @(10, "IndiGo",    "10:40", "02 h 45 m", "₹ 8,983"),
    @(11, "IndiGo",    "13:20", "03 h",      "₹ 8,983"),
    @(12, "IndiGo",    "15:10", "02 h 50 m", "₹ 8,983"),
    @(13, "IndiGo",    "16:35", "02 h 50 m", "₹ 8,983"),
    @(14, "Air India", "16:55", "02 h 55 m", "₹ 8,983"),
    @(15, "Vistara",   "17:15", "02 h 45 m", "₹ 8,983"),
    @(16, "IndiGo",    "18:15", "02 h 50 m", "₹ 8,983"),
    @(17, "IndiGo",    "19:35", "02 h 50 m", "₹ 8,983"),
    @(18, "Air India", "20:15", "02 h 30 m", "₹ 8,983"),
    @(19, "Vistara",   "21:05", "02 h 50 m", "₹ 8,983")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}
